# Update cryptocurrency price/volume data in the sheet to reflect the
# latest scraped values (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.954.17"
$ws.Range("E2").Value = "  -3.70%  "
$ws.Range("D3").Value = "1.637.95"
$ws.Range("E3").Value = "  -6.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9964"
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.73"
$ws.Range("E5").Value = "  -4.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4726"
$ws.Range("E7").Value = "  -6.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2563"
$ws.Range("E8").Value = "  -5.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06014"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07069"
$ws.Range("E10").Value = "  -2.38%  "
$ws.Range("D11").Value = "1.639.76"
$ws.Range("E11").Value = "  -5.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.83"
$ws.Range("E12").Value = "  -1.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6166"
$ws.Range("E13").Value = "  -5.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.373"
$ws.Range("E14").Value = "  -5.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "72.79"
$ws.Range("E15").Value = "  -5.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9977"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "24.946.26"
$ws.Range("E18").Value = "  -3.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006583"
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.18"
$ws.Range("E20").Value = "  -5.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.415"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").Value = "1.845.21"
$ws.Range("E22").Value = "  -6.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.603"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.281"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "133.06"
$ws.Range("E25").Value = "  -2.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.80"
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.361"
$ws.Range("E27").Value = "  -9.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "102.61"
$ws.Range("E28").Value = "  -2.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.660"
$ws.Range("E29").Value = "  -6.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.751"
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07731"
$ws.Range("E31").Value = "  -6.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.558"
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04327"
$ws.Range("E33").Value = "  -7.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9984"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.600"
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9224"
$ws.Range("E36").Value = "  -7.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5824"
$ws.Range("E37").Value = "  -5.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.571"
$ws.Range("E38").Value = "  -6.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01555"
$ws.Range("E39").Value = "  -3.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9979"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8246"
$ws.Range("E41").Value = "  +9.17%  "
$ws.Range("E42").Value = "  -5.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.51"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("E44").Value = "  -4.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.737"
$ws.Range("E45").Value = "  -4.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1107"
$ws.Range("E46").Value = "  -2.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05219"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("E48").Value = "  -3.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.59"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9979"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9997"
$ws.Range("E51").Value = "  -0.46%  "
